$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Cases")

# TestCase_B10 (row 11) Runmode flips from Y to N -- the "Y" marker moves to
# the new TestCase_B12 row (13).
$ws.Range("C11").Value = "N"

# New row 12: TestCase_B11 (copy formatting from row 9, which has the same
# per-column style pattern: s=7,6,7,3)
$ws.Range("A9:D9").Copy($ws.Range("A12:D12"))
$ws.Range("A12").Value = "TestCase_B11"
$ws.Range("B12").Value = "To verify that search,sorting and filtering are retained when user navigates back to search results page from record view page"
$ws.Range("C12").Value = "N"
$ws.Range("D12").Value = "SKIP"
$ws.Rows.Item(12).RowHeight = 16.5

# New row 13: TestCase_B12 (A/C use the s=7 pattern from row 9, B uses the
# s=3 pattern -- matching row 2's B-column style -- and C uses the "Y" value)
$ws.Range("A9").Copy($ws.Range("A13"))
$ws.Range("B2").Copy($ws.Range("B13"))
$ws.Range("C9").Copy($ws.Range("C13"))
$ws.Range("D9").Copy($ws.Range("D13"))
$ws.Range("A13").Value = "TestCase_B12"
$ws.Range("B13").Value = "To verify that the addition of total articles count and total profiles count is equal to total search results count"
$ws.Range("C13").Value = "Y"
$ws.Range("D13").Value = "SKIP"

# Column B width grows and loses its "best fit" flag
$ws.Columns.Item(2).ColumnWidth = 115.140625

# Selection moves
$ws.Range("B6").Select()
